# Auto-generated edit script: refresh market-price-derived profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit* / craft-equivalent columns)
# across all 8 job sheets, per scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1974.6923
$ws.Range("I106").Value = 2107.7273
$ws.Range("K106").Value = 2107.7273
$ws.Range("M106").Value = -1476.7273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 29942.715
$ws.Range("I2").Value = 1065.7826
$ws.Range("J2").Value = 85290.164
$ws.Range("K2").Value = 1065.7826
$ws.Range("L2").Value = 85290.164
$ws.Range("M2").Value = -952.7826
$ws.Range("N2").Value = -85516.164
$ws.Range("H32").Value = 27178.463
$ws.Range("I32").Value = 4564.1636
$ws.Range("K32").Value = 4564.1636
$ws.Range("M32").Value = -4277.1636
$ws.Range("H45").Value = 1488.6052
$ws.Range("I45").Value = 1317.5
$ws.Range("J45").Value = 1967.7
$ws.Range("K45").Value = 1317.5
$ws.Range("L45").Value = 1967.7
$ws.Range("M45").Value = -940.5
$ws.Range("N45").Value = -2721.7
$ws.Range("H55").Value = 9242.857
$ws.Range("J55").Value = 9283.333000000001
$ws.Range("L55").Value = 9283.333000000001
$ws.Range("N55").Value = -9913.333000000001
$ws.Range("H61").Value = 1861.5405
$ws.Range("I61").Value = 962.64703
$ws.Range("K61").Value = 962.64703
$ws.Range("M61").Value = -750.64703
$ws.Range("H74").Value = 2259.7307
$ws.Range("I74").Value = 1178.3
$ws.Range("J74").Value = 2935.625
$ws.Range("K74").Value = 1178.3
$ws.Range("L74").Value = 2935.625
$ws.Range("M74").Value = -304.3
$ws.Range("N74").Value = -4683.625
$ws.Range("H77").Value = 2259.7307
$ws.Range("I77").Value = 1178.3
$ws.Range("J77").Value = 2935.625
$ws.Range("K77").Value = 5891.5
$ws.Range("L77").Value = 14678.125
$ws.Range("M77").Value = -1523.5
$ws.Range("N77").Value = -23414.125
$ws.Range("H80").Value = 14090.571
$ws.Range("I80").Value = 9999
$ws.Range("J80").Value = 14772.5
$ws.Range("K80").Value = 9999
$ws.Range("L80").Value = 14772.5
$ws.Range("M80").Value = -9001
$ws.Range("N80").Value = -16768.5
$ws.Range("H83").Value = 14090.571
$ws.Range("I83").Value = 9999
$ws.Range("J83").Value = 14772.5
$ws.Range("K83").Value = 29997
$ws.Range("L83").Value = 44317.5
$ws.Range("M83").Value = -25005
$ws.Range("N83").Value = -54301.5
$ws.Range("H116").Value = 29942.715
$ws.Range("I116").Value = 1065.7826
$ws.Range("J116").Value = 85290.164
$ws.Range("K116").Value = 1065.7826
$ws.Range("L116").Value = 85290.164
$ws.Range("M116").Value = 1228.2174
$ws.Range("N116").Value = -89878.164
$ws.Range("H136").Value = 1861.5405
$ws.Range("I136").Value = 962.64703
$ws.Range("K136").Value = 2887.94109
$ws.Range("M136").Value = -337.9410899999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 29942.715
$ws.Range("I3").Value = 1065.7826
$ws.Range("J3").Value = 85290.164
$ws.Range("K3").Value = 1065.7826
$ws.Range("L3").Value = 85290.164
$ws.Range("M3").Value = -951.7826
$ws.Range("N3").Value = -85518.164

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 145.88235
$ws.Range("I7").Value = 62.454544
$ws.Range("J7").Value = 298.83334
$ws.Range("K7").Value = 62.454544
$ws.Range("L7").Value = 298.83334
$ws.Range("M7").Value = 50.545456
$ws.Range("N7").Value = -524.83334
$ws.Range("H16").Value = 1221
$ws.Range("I16").Value = 991
$ws.Range("K16").Value = 991
$ws.Range("M16").Value = -704
$ws.Range("H31").Value = 18083.451
$ws.Range("J31").Value = 2310.5854
$ws.Range("L31").Value = 2310.5854
$ws.Range("N31").Value = -2900.5854
$ws.Range("H34").Value = 18083.451
$ws.Range("J34").Value = 2310.5854
$ws.Range("L34").Value = 2310.5854
$ws.Range("N34").Value = -2714.5854
$ws.Range("H113").Value = 1221
$ws.Range("I113").Value = 991
$ws.Range("K113").Value = 991
$ws.Range("M113").Value = 1179
$ws.Range("H132").Value = 2751.2307
$ws.Range("I132").Value = 2733.2632
$ws.Range("J132").Value = 2800
$ws.Range("K132").Value = 8199.7896
$ws.Range("L132").Value = 8400
$ws.Range("M132").Value = -5669.7896
$ws.Range("N132").Value = -13460

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7955.143
$ws.Range("I5").Value = 672.34784
$ws.Range("J5").Value = 41456
$ws.Range("K5").Value = 2017.04352
$ws.Range("L5").Value = 124368
$ws.Range("M5").Value = -1905.04352
$ws.Range("N5").Value = -124592
$ws.Range("H132").Value = 2524.7407
$ws.Range("I132").Value = 2840
$ws.Range("J132").Value = 2453.0908
$ws.Range("K132").Value = 25560
$ws.Range("L132").Value = 22077.8172
$ws.Range("M132").Value = -23030
$ws.Range("N132").Value = -27137.8172
$ws.Range("H135").Value = 7955.143
$ws.Range("I135").Value = 672.34784
$ws.Range("J135").Value = 41456
$ws.Range("K135").Value = 6051.130560000001
$ws.Range("L135").Value = 373104
$ws.Range("M135").Value = -3516.130560000001
$ws.Range("N135").Value = -378174

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 431968.22
$ws.Range("I102").Value = 5180.2856
$ws.Range("J102").Value = 858756.1
$ws.Range("K102").Value = 5180.2856
$ws.Range("L102").Value = 858756.1
$ws.Range("M102").Value = -3558.2856
$ws.Range("N102").Value = -862000.1
$ws.Range("H122").Value = 4230.5713
$ws.Range("I122").Value = 4889.25
$ws.Range("J122").Value = 3352.3333
$ws.Range("K122").Value = 14667.75
$ws.Range("L122").Value = 10056.9999
$ws.Range("M122").Value = -12217.75
$ws.Range("N122").Value = -14956.9999
$ws.Range("H126").Value = 4527060.5
$ws.Range("I126").Value = 2971.6667
$ws.Range("J126").Value = 8404851
$ws.Range("K126").Value = 8915.000100000001
$ws.Range("L126").Value = 25214553
$ws.Range("M126").Value = -6445.000100000001
$ws.Range("N126").Value = -25219493

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2222.2666
$ws.Range("I7").Value = 1642.3334
$ws.Range("K7").Value = 1642.3334
$ws.Range("M7").Value = -1530.3334
$ws.Range("H22").Value = 675
$ws.Range("J22").Value = 675
$ws.Range("L22").Value = 675
$ws.Range("N22").Value = -1265
$ws.Range("H27").Value = 675
$ws.Range("J27").Value = 675
$ws.Range("L27").Value = 675
$ws.Range("N27").Value = -889
$ws.Range("H46").Value = 1267474.9
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1267474.9
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 1267474.9
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -1267850.9
$ws.Range("H126").Value = 2222.2666
$ws.Range("I126").Value = 1642.3334
$ws.Range("K126").Value = 4927.0002
$ws.Range("M126").Value = -2457.0002
$ws.Range("H136").Value = 3536
$ws.Range("I136").Value = 2856
$ws.Range("J136").Value = 4760
$ws.Range("K136").Value = 8568
$ws.Range("L136").Value = 14280
$ws.Range("M136").Value = -6018
$ws.Range("N136").Value = -19380

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 100002110
$ws.Range("I96").Value = 166668880
$ws.Range("J96").Value = 1965.75
$ws.Range("K96").Value = 166668880
$ws.Range("L96").Value = 1965.75
$ws.Range("M96").Value = -166667507
$ws.Range("N96").Value = -4711.75
$ws.Range("H122").Value = 1572
$ws.Range("I122").Value = 1858
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 5574
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -3124
$ws.Range("N122").Value = -7900

Write-Host "Updated 192 cells across $($wb.Worksheets.Count) sheets"
